$wb = $excel.ActiveWorkbook
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}
$values1 = @(3906.399109145206, 0, 48353.76274462014, 0, 289724.0114301849, 9433.134471502228, 0, 2534.277928792104, 0, 0, 0, 0, 0, 2367.37219622158, 1995.762462679798)
for ($i = 0; $i -lt $values1.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value2 = $values1[$i]
}

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}
$values2 = @(6991.052031681918, 0, 197913.7502057619, 0, 289724.0114301849, 16452.51445364119, 0, 8194.52068131253, 0, 0, 0, 0, 0, 7543.193583625169, 6257.586732772244)
for ($i = 0; $i -lt $values2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value2 = $values2[$i]
}

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}
$values3 = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 12889.44547513058, 9263.823477595495)
for ($i = 0; $i -lt $values3.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value2 = $values3[$i]
}

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}
$values4 = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 14047.29522757375, 9263.823477595495)
for ($i = 0; $i -lt $values4.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value2 = $values4[$i]
}

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}
$values5 = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16881.30051551449, 10096.38017359364)
for ($i = 0; $i -lt $values5.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value2 = $values5[$i]
}

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}
$values6 = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16881.30051551449, 10096.38017359364)
for ($i = 0; $i -lt $values6.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value2 = $values6[$i]
}

